# <MS 2/11> Approved (Header Images in Xpath)
#
# - BD_1!I2 gets the postcode value "E126SE" (was blank).
# - The active sheet moves from BD_1 to Intro_0, with the selection on
#   Intro_0 set to D14.

$wb = $excel.ActiveWorkbook

$bd = $wb.Worksheets.Item("BD_1")
$bd.Range("I2").Value = "E126SE"

$intro = $wb.Worksheets.Item("Intro_0")
$intro.Activate()
$intro.Range("D14").Select() | Out-Null
